$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 - copy style (bold, border, centered) from existing
# header cell (H1), then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-12
$values = @{
    2  = @(9, 9)
    3  = @(8, 8)
    4  = @(7, 7)
    5  = @(9, 9)
    6  = @(5, 5)
    7  = @(6, 6)
    8  = @(6, 6)
    9  = @(7, 8)
    10 = @(5, 5)
    11 = @(4, 4)
    12 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
